# Modul button print bill
# Tulis detail item bill (toko, menu, qty, harga) ke worksheet aktif,
# dimulai dari baris 2 (baris 1 berisi header: toko | menu | qty | harga).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$toko = "Kacamata"

$billItems = @(
    @{ menu = "Nasi Hainam Siobak";        qty = 1; harga = 54450.00000000001 },
    @{ menu = "Bakmi Hongkong Siobak";     qty = 1; harga = 52030.00000000001 },
    @{ menu = "Green Tea";                 qty = 1; harga = 18150.0 },
    @{ menu = "Susu Kacang";               qty = 1; harga = 18150.0 }
)

$row = 2
foreach ($item in $billItems) {
    $ws.Cells.Item($row, 1).Value = $toko
    $ws.Cells.Item($row, 2).Value = $item.menu
    $ws.Cells.Item($row, 3).Value = $item.qty
    $ws.Cells.Item($row, 4).Value = $item.harga
    $row++
}
